# fix deceleration problem - change unsigned to signed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 3: new measured/derived constants (signed values) ---
$ws.Range("Q3").Value = 16000
$ws.Range("R3").Value = 3388705
$ws.Range("S3").Value = 564784
$ws.Range("T3").Value = 8094114

# --- R6 gets a comma/number style so the (now much larger) squared value is readable ---
$ws.Range("R6").Style = "Komma"

# --- R10:R20 frqO formulas: add the missing "* 1000" factor (unsigned -> signed fix) ---
$ws.Range("R10").Formula = "=SQRT(`$Q`$6 * Q10 * 1000 + `$R`$6)"
$ws.Range("R11").Formula = "=SQRT(`$Q`$6 * Q11 * 1000 + `$R`$6)"
$ws.Range("R12").Formula = "=SQRT(`$Q`$6 * Q12 * 1000 + `$R`$6)"
$ws.Range("R13").Formula = "=SQRT(`$Q`$6 * Q13 * 1000 + `$R`$6)"
$ws.Range("R14").Formula = "=SQRT(`$Q`$6 * Q14 * 1000 + `$R`$6)"
$ws.Range("R15").Formula = "=SQRT(`$Q`$6 * Q15 * 1000 + `$R`$6)"
$ws.Range("R16").Formula = "=SQRT(`$Q`$6 * Q16 * 1000 + `$R`$6)"
$ws.Range("R17").Formula = "=SQRT(`$Q`$6 * Q17 * 1000 + `$R`$6)"
$ws.Range("R18").Formula = "=SQRT(`$Q`$6 * Q18 * 1000 + `$R`$6)"
$ws.Range("R19").Formula = "=SQRT(`$Q`$6 * Q19 * 1000 + `$R`$6)"
$ws.Range("R20").Formula = "=SQRT(`$Q`$6 * Q20 * 1000 + `$R`$6)"

# --- Column widths: widen R (frqO) and give S (frqA2) its own explicit width ---
$ws.Columns("R").ColumnWidth = 23
$ws.Columns("S").ColumnWidth = 16.28515625

# --- View: scroll back to top-left and move selection to T3 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("T3").Select()
